# Correcion entrada index entero
# The "HomeWork 1" header (cell B1) had a space; replace it with an
# underscore so it reads "HomeWork_1".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "HomeWork_1"

# Leave the selection on B2, matching where the user ended up after
# editing the header cell.
$ws.Range("B2").Select()
